$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 440, shifting existing rows 440:514 down to 441:515
$ws.Rows("440:440").Insert()

# Populate the new row 440 with the new weekly data point
$ws.Range("A440").Value = 4
$ws.Range("B440").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C440").Value = "Los Lagos"
$ws.Range("D440").Value = 45180
$ws.Range("E440").Value = 10
$ws.Range("F440").Value = 100112037
$ws.Range("G440").Value = "Cebollín"
$ws.Range("H440").Value = "Sin especificar"
$ws.Range("I440").Value = "Primera"
$ws.Range("J440").Value = 70
$ws.Range("K440").Value = 7000
$ws.Range("L440").Value = 7000
$ws.Range("M440").Value = 7000
$ws.Range("N440").Value = "$/paquete 36 unidades"
$ws.Range("O440").Value = "Región Metropolitana"
$ws.Range("P440").Value = 194
$ws.Range("Q440").Value = 36
$ws.Range("R440").Value = "Hortaliza"
